$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Job Type" column header
$ws.Range("E4").Value = "Job Type"
$ws.Range("E4").Font.Bold = $true

# Fill in Job Type values for existing rows (5-11)
$ws.Range("E5").Value = "Social"
$ws.Range("E6").Value = "Functional"
$ws.Range("E7").Value = "Functional"
$ws.Range("E8").Value = "Functional"
$ws.Range("E9").Value = "Social"
$ws.Range("E10").Value = "Emotional"
$ws.Range("E11").Value = "Functional"

# Update D10 text
$ws.Range("D10").Value = "for leisure activities I like"

# Add new row 12
$ws.Range("B12").Value = "Energy"
$ws.Range("C12").Value = "for exercise after the day ends"
$ws.Range("D12").Value = "so I can get in shape"
$ws.Range("E12").Value = "Emotional"

# Auto-fit columns C and D to match the widened content
$ws.Columns("C:D").AutoFit()

# Update selection to match target state
$ws.Range("E19").Select()
